# Weekly update: insert a new week's worth of price rows (239-244) for
# "Melón" at "Vega Monumental Concepción", pushing the previous week's
# rows (old 239-244) down to 245-250.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows at position 239; this shifts the old rows
# 239-244 down to 245-250 and extends the sheet dimension accordingly.
$ws.Rows("239:244").Insert()

# New week's data (Fecha serial 44615) for the 6 Melón variety/quality
# combinations, mirroring the rows that were pushed down.

# Row 239: Calameño - Extra
$ws.Range("A239").Value = 11
$ws.Range("B239").Value = "Vega Monumental Concepción"
$ws.Range("C239").Value = "Bíobío"
$ws.Range("D239").Value = 44615
$ws.Range("E239").Value = 8
$ws.Range("F239").Value = 100112027
$ws.Range("G239").Value = "Melón"
$ws.Range("H239").Value = "Calameño"
$ws.Range("I239").Value = "Extra"
$ws.Range("J239").Value = 1000
$ws.Range("K239").Value = 1000
$ws.Range("L239").Value = 1000
$ws.Range("M239").Value = 1000
$ws.Range("N239").Value = "$/unidad"
$ws.Range("O239").Value = "Región de O'Higgins"
$ws.Range("P239").Value = 1000
$ws.Range("Q239").Value = 1
$ws.Range("R239").Value = "Hortaliza"

# Row 240: Calameño - Primera
$ws.Range("A240").Value = 11
$ws.Range("B240").Value = "Vega Monumental Concepción"
$ws.Range("C240").Value = "Bíobío"
$ws.Range("D240").Value = 44615
$ws.Range("E240").Value = 8
$ws.Range("F240").Value = 100112027
$ws.Range("G240").Value = "Melón"
$ws.Range("H240").Value = "Calameño"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 1500
$ws.Range("K240").Value = 800
$ws.Range("L240").Value = 800
$ws.Range("M240").Value = 800
$ws.Range("N240").Value = "$/unidad"
$ws.Range("O240").Value = "Región de O'Higgins"
$ws.Range("P240").Value = 800
$ws.Range("Q240").Value = 1
$ws.Range("R240").Value = "Hortaliza"

# Row 241: Calameño - Segunda
$ws.Range("A241").Value = 11
$ws.Range("B241").Value = "Vega Monumental Concepción"
$ws.Range("C241").Value = "Bíobío"
$ws.Range("D241").Value = 44615
$ws.Range("E241").Value = 8
$ws.Range("F241").Value = 100112027
$ws.Range("G241").Value = "Melón"
$ws.Range("H241").Value = "Calameño"
$ws.Range("I241").Value = "Segunda"
$ws.Range("J241").Value = 1500
$ws.Range("K241").Value = 700
$ws.Range("L241").Value = 700
$ws.Range("M241").Value = 700
$ws.Range("N241").Value = "$/unidad"
$ws.Range("O241").Value = "Región de O'Higgins"
$ws.Range("P241").Value = 700
$ws.Range("Q241").Value = 1
$ws.Range("R241").Value = "Hortaliza"

# Row 242: Tuna - Extra
$ws.Range("A242").Value = 11
$ws.Range("B242").Value = "Vega Monumental Concepción"
$ws.Range("C242").Value = "Bíobío"
$ws.Range("D242").Value = 44615
$ws.Range("E242").Value = 8
$ws.Range("F242").Value = 100112027
$ws.Range("G242").Value = "Melón"
$ws.Range("H242").Value = "Tuna"
$ws.Range("I242").Value = "Extra"
$ws.Range("J242").Value = 1000
$ws.Range("K242").Value = 1000
$ws.Range("L242").Value = 1000
$ws.Range("M242").Value = 1000
$ws.Range("N242").Value = "$/unidad"
$ws.Range("O242").Value = "Región de O'Higgins"
$ws.Range("P242").Value = 1000
$ws.Range("Q242").Value = 1
$ws.Range("R242").Value = "Hortaliza"

# Row 243: Tuna - Primera
$ws.Range("A243").Value = 11
$ws.Range("B243").Value = "Vega Monumental Concepción"
$ws.Range("C243").Value = "Bíobío"
$ws.Range("D243").Value = 44615
$ws.Range("E243").Value = 8
$ws.Range("F243").Value = 100112027
$ws.Range("G243").Value = "Melón"
$ws.Range("H243").Value = "Tuna"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 1500
$ws.Range("K243").Value = 800
$ws.Range("L243").Value = 800
$ws.Range("M243").Value = 800
$ws.Range("N243").Value = "$/unidad"
$ws.Range("O243").Value = "Región de O'Higgins"
$ws.Range("P243").Value = 800
$ws.Range("Q243").Value = 1
$ws.Range("R243").Value = "Hortaliza"

# Row 244: Tuna - Segunda
$ws.Range("A244").Value = 11
$ws.Range("B244").Value = "Vega Monumental Concepción"
$ws.Range("C244").Value = "Bíobío"
$ws.Range("D244").Value = 44615
$ws.Range("E244").Value = 8
$ws.Range("F244").Value = 100112027
$ws.Range("G244").Value = "Melón"
$ws.Range("H244").Value = "Tuna"
$ws.Range("I244").Value = "Segunda"
$ws.Range("J244").Value = 1500
$ws.Range("K244").Value = 700
$ws.Range("L244").Value = 700
$ws.Range("M244").Value = 700
$ws.Range("N244").Value = "$/unidad"
$ws.Range("O244").Value = "Región de O'Higgins"
$ws.Range("P244").Value = 700
$ws.Range("Q244").Value = 1
$ws.Range("R244").Value = "Hortaliza"
